$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text format for price (D) and volume/percentage (E) columns so that
# numeric-looking strings (e.g. "0.001500", "5.42%") are preserved exactly as
# text rather than being re-interpreted as numbers/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.42%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.26%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.216"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.77%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.35%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.378"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "7.24%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.051"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.35%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9336"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.43%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1009"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "10.01%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1843"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9.09%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08739"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.03%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03316"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.32%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09899"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.42%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001500"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005720"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.78%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.470"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.50%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.970"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "6.39%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.151"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.78%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3378"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.38%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1302"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.31%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.304"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.05%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2228"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.00%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04575"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.83%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.72%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004429"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.61%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01771"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "13.32%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.48%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007737"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.08%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1410"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.32%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007140"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-24.82%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002280"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.62%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009192"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.90%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005994"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.00%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "16.48%"
